$wb = $excel.ActiveWorkbook

$wsRisiken = $wb.Worksheets.Item("Risiken")
$wsRisiken.Range("C6").Value = "Kinect erkennt Menschen nicht, die sich parallel zur Wand ausgerichtet vor der Video Wall bewegen."
$wsRisiken.Range("B6").Value = "Kinect: Erkennung von der Seite"

$wsHistory = $wb.Worksheets.Item("Änderungsgeschichte")
$wsHistory.Range("C9").Value = 'Risiko 3: "Kinect: Erkennung von der Seite" ist bereinigt durch die Aufnahmen, die mit Kinect im Gebäude 4 durchgeführt wurden (siehe Dokument Vorstudie).'
$wsHistory.Rows.Item(9).RowHeight = 45

$wsHistory.Range("C9").Select()
$wsRisiken.Range("B6").Select()
